$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price observation was inserted before the existing row 289, which
# pushes every following row (old 289..359) down by one (new 290..360).
$ws.Rows.Item(289).Insert()

# Seed the new row 289 from the row that landed right below it (old row 289,
# now at 290) so the static descriptive columns (market/region/category/...)
# carry over, then overwrite the columns that hold the new observation's
# own data (date, volume, min/max/avg price, price per kg).
$ws.Range("A289:R289").Value2 = $ws.Range("A290:R290").Value2

$ws.Range("D289").Value2 = 44782
$ws.Range("J289").Value2 = 20
$ws.Range("K289").Value2 = 10000
$ws.Range("L289").Value2 = 10000
$ws.Range("M289").Value2 = 10000
$ws.Range("P289").Value2 = 833

# Match the date cell's number format used by the other date cells.
$ws.Range("D289").NumberFormat = $ws.Range("D290").NumberFormat
